$d = $word.ActiveDocument

# The "suivi personnel" table is the second (last) table in the document.
$t = $d.Tables.Item(2)

# --- 1) Append a new row for 20/02/2018 ------------------------------------
# Rows.Add() clones the formatting (pStyle/numPr) of the row above it, which
# is exactly what the second column of this table needs.
$newRow = $t.Rows.Add()
$newRow.Cells.Item(1).Range.Text = "20/02/2018"

$cell2 = $newRow.Cells.Item(2)
$para2 = $cell2.Range.Paragraphs.Item(1)
# Append a temporary trailing marker character; it lets us later add a
# collapsed bookmark right after the real text without hitting the
# paragraph-end / cell-end boundary (removed again in step 3).
$para2.Range.Text = "Elaboration du fichier communZ"

# --- 2) Remove the existing "_GoBack" bookmark from the 19/02/2018 row -----
# It currently sits between ", ouvert/fermé " and "(vu avec M. Dugast)".
# Bookmark.Delete() is a no-op in this host, but overwriting a range whose
# interior strictly contains the (zero-width) bookmark point does clear it,
# so round-trip two characters straddling the bookmark through a dummy value.
$bm = $d.Bookmarks.Item("_GoBack")
$bmStart = $bm.Start
$straddle = $d.Range($bmStart - 1, $bmStart + 1)
$savedText = $straddle.Text
$straddle.Text = "ZZ"
$restore = $d.Range($bmStart - 1, $bmStart + 1)
$restore.Text = $savedText

# --- 3) Add the "_GoBack" bookmark at the end of the new row's text -------
# Re-fetch the paragraph fresh since positions shifted after step 2.
$para2Fresh = $t.Rows.Item(3).Cells.Item(2).Range.Paragraphs.Item(1)
$markerPos = $para2Fresh.Range.End - 2   # just after "commun", before "Z"
$bmRange = $d.Range($markerPos, $markerPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

# Remove the temporary trailing "Z" marker character.
$para2Final = $t.Rows.Item(3).Cells.Item(2).Range.Paragraphs.Item(1)
$zEnd = $para2Final.Range.End - 1
$zRange = $d.Range($zEnd - 1, $zEnd)
$zRange.Delete()
